$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$q1 = [char]0x201C
$q2 = [char]0x201D

# Locate the paragraph that currently holds the "Submit the metric..." sentence.
$idx = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Submit the metric and your code*") {
        $idx = $i
        break
    }
    $i = $i + 1
}

# 1) Turn that paragraph into the new bold "Submission:" heading paragraph.
$target = $d.Paragraphs.Item($idx)
$xmlSubmission = "<w:p $wNs>" + `
    "<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>" + `
    "<w:r><w:rPr><w:b/></w:rPr><w:t>Submission</w:t></w:r>" + `
    "<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>" + `
    "</w:p>"
$target.Range.InsertXML($xmlSubmission)

# 2) Insert the new explanatory paragraph right after it (plain, no explicit spacing pPr).
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$xmlInformative = "<w:p $wNs>" + `
    "<w:r><w:t>The metric value should be computed for each elapsed time step (by calling the provided code or by implementing yourself). The metric value should be reported for several elapsed time steps. The number of elapsed time steps should be sufficient to establish an " + $q1 + "informative profile" + $q2 + ".</w:t></w:r>" + `
    "</w:p>"
$target.Range.InsertXML($xmlInformative)

# 3) Insert the "For further details..." paragraph.
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$xmlFurther = "<w:p $wNs>" + `
    "<w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr>" + `
    "<w:r><w:t>For further details regarding submission of the metric and your code, please refer to the main CP4 problem description document, e.g. PPAML-Challenge-Problem-4.pdf.</w:t></w:r>" + `
    "</w:p>"
$target.Range.InsertXML($xmlFurther)

# 4) Blank spacer paragraph.
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$xmlBlank1 = "<w:p $wNs><w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr></w:p>"
$target.Range.InsertXML($xmlBlank1)

# 5) "Sample output files..." paragraph with proofErr spellcheck markers.
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$xmlSample = "<w:p $wNs>" + `
    "<w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr>" + `
    "<w:r><w:t>Sample output files for this problem have been provided in the " + $q1 + "</w:t></w:r>" + `
    "<w:proofErr w:type=`"spellStart`"/>" + `
    "<w:r><w:t>sampleoutput</w:t></w:r>" + `
    "<w:proofErr w:type=`"spellEnd`"/>" + `
    "<w:r><w:t>" + $q2 + " folder:</w:t></w:r>" + `
    "</w:p>"
$target.Range.InsertXML($xmlSample)

# 6) Blank spacer paragraph.
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$xmlBlank2 = "<w:p $wNs><w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr></w:p>"
$target.Range.InsertXML($xmlBlank2)

# 7) The Courier-New filename paragraph with proofErr gramStart/gramEnd markers.
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$rFonts = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'
$xmlCsv = "<w:p $wNs>" + `
    "<w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/><w:rPr>$rFonts</w:rPr></w:pPr>" + `
    "<w:proofErr w:type=`"gramStart`"/>" + `
    "<w:r><w:rPr>$rFonts</w:rPr><w:t>problem-</w:t></w:r>" + `
    "<w:r><w:rPr>$rFonts</w:rPr><w:t>7</w:t></w:r>" + `
    "<w:r><w:rPr>$rFonts</w:rPr><w:t>-query-1-metric-1.csv</w:t></w:r>" + `
    "<w:proofErr w:type=`"gramEnd`"/>" + `
    "</w:p>"
$target.Range.InsertXML($xmlCsv)

# 8) Blank spacer paragraph.
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$idx = $idx + 1
$target = $d.Paragraphs.Item($idx)
$xmlBlank3 = "<w:p $wNs><w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr></w:p>"
$target.Range.InsertXML($xmlBlank3)

# 9) Remove the two now-orphaned legacy blank paragraphs (the bottom-bordered one and the
#    bold empty one) that used to separate "Submit the metric..." from "Ground Truth Details:".
$borderPara = $d.Paragraphs.Item($idx + 1)
$borderPara.Range.Delete()
$boldPara = $d.Paragraphs.Item($idx + 1)
$boldPara.Range.Delete()

# 10) Re-home the _GoBack bookmark onto the start of the "Ground Truth Details:" paragraph.
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}
$gtd = $d.Paragraphs.Item($idx + 1)
Write-Host "gtd text: [$($gtd.Range.Text)]"
$d.Bookmarks.Add("_GoBack", $d.Range($gtd.Range.Start, $gtd.Range.Start))

Write-Host "done"
